$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1 "Fecha", matching the style of the other header cells (e.g. J1)
$ws.Cells.Item(1, 10).Copy()
$ws.Cells.Item(1, 11).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 11).Value = "Fecha"

# New data column K2:K23 "09/08/24" as plain text (avoid Excel auto-converting it to a date)
$dataRange = $ws.Range("K2:K23")
$dataRange.NumberFormat = "@"
$dataRange.Value = "09/08/24"
$dataRange.Style = "Normal"
